# Add a new school record ("FUNDAÇÃO GETULIO VARGAS") right above the
# existing "CENTRO DE ENSINO MEDIO ARY RIBEIRO VALADAO FILHO" row (row 98),
# pushing that row and everything below it down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 98 (shifts old rows 98:151 down to 99:152).
$ws.Range("A98").EntireRow.Insert()

# Populate the newly inserted row with the new record.
$ws.Range("A98").Value = "FUNDAÇÃO GETULIO VARGAS"
$ws.Range("B98").Value = "123.456.789-10"
$ws.Range("D98").Value = 123456

# Update the view state to match where the author ended up after editing.
$ws.Range("D98").Select()
